# Update "想去人数" (want-to-go count) figures pulled from the latest
# bilibili scrape (gh-pages output regenerated at commit 456a3b4).
#
# Sheet "展览" (Exhibition):
#   F6  2528 -> 2531
#   F14  393 -> 394
#   F16  492 -> 493
#
# Sheet "本地生活" (Local life):
#   F2  6073 -> 6074
#   F5  1374 -> 1376
#   F6    34 -> 36
#
# Sheet "全部类型" (All types, union of the other sheets):
#   F2  6073 -> 6074
#   F5  1374 -> 1376
#   F14 2528 -> 2531
#   F17   34 -> 36
#   F27  393 -> 394
#   F31  492 -> 493

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F6").Value = 2531
$wsExhibition.Range("F14").Value = 394
$wsExhibition.Range("F16").Value = 493

$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 6074
$wsLocalLife.Range("F5").Value = 1376
$wsLocalLife.Range("F6").Value = 36

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6074
$wsAll.Range("F5").Value = 1376
$wsAll.Range("F14").Value = 2531
$wsAll.Range("F17").Value = 36
$wsAll.Range("F27").Value = 394
$wsAll.Range("F31").Value = 493
